# Refresh the "cryptos" list: updates the Price (column D) and Volume(1h)
# (column E) columns for the rows whose figures moved since the last
# GitHub Actions run.
#
# Percentage cells in column E are always text already (leading/trailing
# spaces + a trailing "%"), so a plain assignment keeps them text.
#
# Some new Price values in column D are plain decimals (e.g. "7.20",
# "2.20") that Excel's COM layer would happily reinterpret as numbers -
# silently dropping the trailing zero / turning the cell numeric. To keep
# them as literal text (matching the original inline-string cells), those
# are written with a leading single-quote (Excel's "force text" marker),
# and the cell style is then reset to "Normal" so no stray formatting is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "48.088.86"; E = "  +1.36%  " },
    @{ Row = 3;  D = "2.510.88";  E = "  +0.04%  " },
    @{ Row = 4;                   E = "  +0.00%  " },
    @{ Row = 5;  D = "323.28";    E = "  -0.34%  " },
    @{ Row = 6;  D = "108.74";    E = "  -0.79%  " },
    @{ Row = 7;  D = "0.525";     E = "  -0.53%  " },
    @{ Row = 9;                   E = "  +3.49%  " },
    @{ Row = 10; D = "40.35";     E = "  +2.92%  " },
    @{ Row = 11; D = "19.68";     E = "  +5.34%  " },
    @{ Row = 12; D = "0.0816";    E = "  -0.46%  " },
    @{ Row = 13;                  E = "  +0.67%  " },
    @{ Row = 14; D = "7.20";      E = "  -0.40%  " },
    @{ Row = 15; D = "2.901.56";  E = "  +0.03%  " },
    @{ Row = 16; D = "2.513.60";  E = "  +0.25%  " },
    @{ Row = 17; D = "0.852";     E = "  -1.83%  " },
    @{ Row = 18; D = "47.936.77"; E = "  +1.14%  " },
    @{ Row = 19; D = "13.38";     E = "  +3.07%  " },
    @{ Row = 20; D = "6.62";      E = "  -1.65%  " },
    @{ Row = 21;                  E = "  -0.73%  " },
    @{ Row = 22;                  E = "  +4.67%  " },
    @{ Row = 23;                  E = "  -0.13%  " },
    @{ Row = 24; D = "247.93";    E = "  -0.87%  " },
    @{ Row = 25;                  E = "  -1.18%  " },
    @{ Row = 26;                  E = "  +0.04%  " },
    @{ Row = 27; D = "25.93";     E = "  -1.22%  " },
    @{ Row = 28; D = "10.24";     E = "  +1.61%  " },
    @{ Row = 29; D = "2.20";      E = "  -4.27%  " },
    @{ Row = 30;                  E = "  +4.68%  " },
    @{ Row = 31; D = "35.24";     E = "  -2.14%  " },
    @{ Row = 32; D = "49.84";     E = "  -0.92%  " },
    @{ Row = 33;                  E = "  -0.04%  " },
    @{ Row = 34; D = "5.39";      E = "  -1.28%  " },
    @{ Row = 35;                  E = "  +0.06%  " },
    @{ Row = 36; D = "0.0786";    E = "  -1.25%  " },
    @{ Row = 37;                  E = "  -1.49%  " },
    @{ Row = 38;                  E = "  -1.40%  " },
    @{ Row = 39;                  E = "  -1.29%  " },
    @{ Row = 40;                  E = "  -0.59%  " },
    @{ Row = 41; D = "22.40";     E = "  +4.14%  " },
    @{ Row = 42; D = "118.58";    E = "  -4.06%  " },
    @{ Row = 43; D = "2.18";      E = "  -3.24%  " },
    @{ Row = 44; D = "0.0298";    E = "  -0.32%  " },
    @{ Row = 45; D = "1.997.67";  E = "  -0.02%  " },
    @{ Row = 46; D = "3.13";      E = "  +0.89%  " },
    @{ Row = 47;                  E = "  -2.94%  " },
    @{ Row = 48; D = "1.82";      E = "  +1.43%  " },
    @{ Row = 49; D = "9.11";      E = "  +0.17%  " },
    @{ Row = 50; D = "5.20";      E = "  -1.89%  " },
    @{ Row = 51; D = "56.56";     E = "  +0.75%  " }
)

# NOTE: the value-setting logic is inlined in the loop body (rather than
# calling a helper function) because this COM host's PowerShell engine
# mishandles a user-defined function invoked more than once from inside a
# `foreach` block.
foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey("D")) {
        $addr = "D" + $row
        $val = $u.D
        if ($val -match '^[+-]?\d+(\.\d+)?$') {
            # Plain-decimal text: guard against COM auto-coercing the
            # string into a Number (which would drop trailing zeros).
            $ws.Range($addr).Value = "'" + $val
            $ws.Range($addr).Style = "Normal"
        } else {
            $ws.Range($addr).Value = $val
        }
    }

    if ($u.ContainsKey("E")) {
        $ws.Range("E" + $row).Value = $u.E
    }
}
